$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.320.60'
$ws.Range("D3").Value = '1.690.94'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '218.51'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").Value = '0.5270'
$ws.Range("E6").Value = '  +4.40%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.2709'
$ws.Range("E8").Value = '  +2.19%  '
$ws.Range("D9").Value = '0.06429'
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").Value = '22.08'
$ws.Range("E10").Value = '  +2.84%  '
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").Value = '1.714.65'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("D13").Value = '4.572'
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").Value = '0.5858'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").Value = '0.000008517'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = '64.54'
$ws.Range("D17").Value = '26.381.32'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").Value = '4.948'
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D20").Value = '10.90'
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("D21").Value = '189.67'
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").Value = '6.220'
$ws.Range("E22").Value = '  +0.65%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '144.80'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").Value = '7.707'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("E26").Value = '  +5.41%  '
$ws.Range("D27").Value = '15.86'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("D28").Value = '0.06728'
$ws.Range("E28").Value = '  +16.15%  '
$ws.Range("D29").Value = '1.357'
$ws.Range("E29").Value = '  +6.24%  '
$ws.Range("D30").Value = '1.330'
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("D31").Value = '3.589'
$ws.Range("E31").Value = '  +2.31%  '
$ws.Range("D32").Value = '3.575'
$ws.Range("E32").Value = '  +1.38%  '
$ws.Range("E33").Value = '  +2.57%  '
$ws.Range("D34").Value = '1.031'
$ws.Range("E34").Value = '  +2.28%  '
$ws.Range("D35").Value = '0.6244'
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("D36").Value = '2.396'
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("E37").Value = '  +2.35%  '
$ws.Range("D38").Value = '6.371'
$ws.Range("E38").Value = '  +5.97%  '
$ws.Range("D39").Value = '1.116.18'
$ws.Range("E39").Value = '  +4.30%  '
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").Value = '0.8902'
$ws.Range("E41").Value = '  +3.53%  '
$ws.Range("D42").Value = '1.017'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").Value = '100.87'
$ws.Range("E43").Value = '  +1.50%  '
$ws.Range("D44").Value = '1.839.53'
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("E45").Value = '  +3.87%  '
$ws.Range("D46").Value = '57.02'
$ws.Range("E46").Value = '  +2.65%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '8.206'
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("D49").Value = '0.05270'
$ws.Range("E49").Value = '  +1.73%  '
$ws.Range("D50").Value = '0.4305'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '6.084'
$ws.Range("E51").Value = '  +3.98%  '
